$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" column in H, matching the header style used by the
# other headers (B1:G1).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Fill in the per-row Save values.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 0
